# Apply updated violent-crime YTD figures for 2022-12-17
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 116
$ws.Range("D3").Value = 139
$ws.Range("E3").Value = 148
$ws.Range("H3").Value = 158
$ws.Range("I3").Value = 196
$ws.Range("H4").Value = 14
$ws.Range("B6").Value = 384
$ws.Range("C6").Value = 491
$ws.Range("D6").Value = 426
$ws.Range("E6").Value = 489
$ws.Range("F6").Value = 557
$ws.Range("G6").Value = 440
$ws.Range("H6").Value = 453
$ws.Range("I6").Value = 508
$ws.Range("B7").Value = 518
$ws.Range("C7").Value = 648
$ws.Range("D7").Value = 666
$ws.Range("E7").Value = 723
$ws.Range("F7").Value = 805
$ws.Range("G7").Value = 674
$ws.Range("H7").Value = 737
$ws.Range("I7").Value = 846

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("D3").Value = 9
$ws.Range("E6").Value = 54
$ws.Range("F6").Value = 52
$ws.Range("H6").Value = 38
$ws.Range("D7").Value = 49
$ws.Range("E7").Value = 67
$ws.Range("F7").Value = 61
$ws.Range("H7").Value = 49

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I6").Value = 21
$ws.Range("I7").Value = 37

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("D5").Value = 11
$ws.Range("D6").Value = 14

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("E4").Value = 6
$ws.Range("E5").Value = 10

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("H4").Value = 2
$ws.Range("B6").Value = 32
$ws.Range("F6").Value = 39
$ws.Range("B7").Value = 37
$ws.Range("F7").Value = 59
$ws.Range("H7").Value = 47

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("E8").Value = 54
$ws.Range("I8").Value = 41
$ws.Range("G10").Value = 5
$ws.Range("I19").Value = 23
$ws.Range("B28").Value = 37
$ws.Range("F28").Value = 59
$ws.Range("H28").Value = 47
$ws.Range("D32").Value = 49
$ws.Range("E32").Value = 67
$ws.Range("F32").Value = 61
$ws.Range("H32").Value = 49
$ws.Range("I36").Value = 37
$ws.Range("H47").Value = 25
$ws.Range("I47").Value = 26
$ws.Range("D50").Value = 14
$ws.Range("C51").Value = 3
$ws.Range("D53").Value = 75
$ws.Range("E53").Value = 88
$ws.Range("H53").Value = 106
$ws.Range("I53").Value = 125
$ws.Range("D65").Value = 28
$ws.Range("F65").Value = 40
$ws.Range("H70").Value = 17
$ws.Range("F77").Value = 23
$ws.Range("D85").Value = 6
$ws.Range("E88").Value = 10
$ws.Range("C94").Value = 6
$ws.Range("B98").Value = 518
$ws.Range("C98").Value = 648
$ws.Range("D98").Value = 666
$ws.Range("E98").Value = 723
$ws.Range("F98").Value = 805
$ws.Range("G98").Value = 674
$ws.Range("H98").Value = 737
$ws.Range("I98").Value = 846

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 14
$ws.Range("E3").Value = 17
$ws.Range("H3").Value = 22
$ws.Range("D6").Value = 46
$ws.Range("H6").Value = 68
$ws.Range("D7").Value = 75
$ws.Range("E7").Value = 88
$ws.Range("H7").Value = 106
$ws.Range("I7").Value = 125

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("D5").Value = 27
$ws.Range("F5").Value = 33
$ws.Range("D6").Value = 28
$ws.Range("F6").Value = 40

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("F6").Value = 12
$ws.Range("F7").Value = 23

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 5

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I6").Value = 15
$ws.Range("I7").Value = 23

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 3

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("D4").Value = 4
$ws.Range("D5").Value = 6

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I3").Value = 2
$ws.Range("H5").Value = 15
$ws.Range("H6").Value = 25
$ws.Range("I6").Value = 26

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("H4").Value = 12
$ws.Range("H5").Value = 17

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("C4").Value = 5
$ws.Range("C5").Value = 6

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("E5").Value = 42
$ws.Range("I5").Value = 31
$ws.Range("E6").Value = 54
$ws.Range("I6").Value = 41
